$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

# Snapshot all source rows (2-31) before any writes, since this edit permutes rows
$snapshot = @{}
for ($r = 2; $r -le 31; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2()
    }
    $snapshot[$r] = $rowVals
}

# Target row -> source row mapping (the weekly refresh reshuffles existing records)
$perm = @{
    2 = 4
    3 = 10
    4 = 19
    5 = 30
    6 = 16
    7 = 15
    8 = 6
    9 = 23
    10 = 9
    11 = 26
    12 = 17
    13 = 13
    14 = 22
    15 = 25
    16 = 5
    17 = 12
    18 = 21
    19 = 18
    20 = 14
    21 = 20
    22 = 29
    23 = 27
    24 = 11
    25 = 28
    26 = 31
    27 = 2
    28 = 3
    29 = 24
    30 = 7
    31 = 8
}

foreach ($target in ($perm.Keys | Sort-Object)) {
    $source = $perm[$target]
    $srcRow = $snapshot[$source]
    foreach ($col in $cols) {
        $ws.Range("$col$target").Value = $srcRow[$col]
    }
}

